# Apply crypto price/volume updates as described in the commit diff.
# Values in column D that look numeric (e.g. "1.005") must be written
# with a leading apostrophe so Excel keeps them as text (matching the
# original inline-string cell type) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.838.87'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '1.904.83'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'313.04"
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = "'0.5033"
$ws.Range("E7").Value = '  +4.08%  '
$ws.Range("D8").Value = "'0.3818"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").Value = "'0.9078"
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("D11").Value = "'20.78"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.951.18'
$ws.Range("E12").Value = '  +1.90%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = "'0.07652"
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("D14").Value = "'5.474"
$ws.Range("E14").Value = '  -0.74%  '
$ws.Range("D15").Value = "'6.586"
$ws.Range("E15").Value = '  -0.82%  '
$ws.Range("D16").Value = "'91.33"
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = "'0.000008712"
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").Value = '27.866.68'
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").Value = "'14.51"
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").Value = "'5.167"
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").Value = "'10.81"
$ws.Range("E23").Value = '  -0.87%  '
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("E25").Value = '  -2.31%  '
$ws.Range("D26").Value = "'2.224"
$ws.Range("E26").Value = '  +5.14%  '
$ws.Range("D27").Value = "'18.36"
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("D28").Value = "'115.21"
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").Value = "'4.902"
$ws.Range("E29").Value = '  -0.96%  '
$ws.Range("D30").Value = "'0.08985"
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("D31").Value = "'3.203"
$ws.Range("E31").Value = '  -4.14%  '
$ws.Range("D32").Value = "'1.230"
$ws.Range("E32").Value = '  -1.84%  '
$ws.Range("D33").Value = "'0.7630"
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("D34").Value = "'4.635"
$ws.Range("E34").Value = '  -1.02%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").Value = "'2.536"
$ws.Range("E36").Value = '  -3.21%  '
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("D38").Value = "'0.5556"
$ws.Range("E38").Value = '  +1.01%  '
$ws.Range("D39").Value = "'3.018"
$ws.Range("E39").Value = '  +1.20%  '
$ws.Range("D40").Value = "'0.05253"
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("D41").Value = "'6.968"
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("D42").Value = "'8.465"
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("E43").Value = '  -0.70%  '
$ws.Range("D44").Value = "'111.42"
$ws.Range("E44").Value = '  +3.74%  '
$ws.Range("D45").Value = "'10.60"
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("D46").Value = "'0.4790"
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("D48").Value = "'1.629"
$ws.Range("E48").Value = '  -1.73%  '
$ws.Range("D49").Value = "'67.26"
$ws.Range("E49").Value = '  -1.64%  '
$ws.Range("D50").Value = "'0.06077"
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").Value = "'0.8992"
$ws.Range("E51").Value = '  -0.23%  '
